{"js": "const body = context.document.body;\n\n// Collapsed range at the very end of the body (right after \"git hub\",\n// before the body's trailing paragraph mark).\nconst endRange = body.getRange(\"End\");\n\n// Flat-OPC wrapped WordprocessingML: an empty paragraph that keeps the\n// \"en-US\" paragraph-mark language already in use, followed by a fresh\n// paragraph with the new second line of text. insertOoxml(..., \"After\")\n// reproduces this exact shape - no stray empty run and no inherited\n// formatting leaking into the new text paragraph.\nconst xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p><w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr></w:p>' +\n  '<w:p><w:r><w:t>\u0414\u043e\u0431\u0430\u0432\u043b\u044f\u044e \u0432\u0442\u043e\u0440\u0443\u044e \u0441\u0442\u0440\u043e\u0447\u043a\u0443.</w:t></w:r></w:p>' +\n  '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\nendRange.insertOoxml(xml, \"After\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Collapsed range sitting exactly at the end of the document body (after\n# \"\u0424\u0430\u0439\u043b \u0434\u043b\u044f git hub\", before the final paragraph mark).\n$endPos = $d.Content.End\n$r = $d.Range($endPos, $endPos)\n\n# Flat-OPC wrapped WordprocessingML fragment: an empty paragraph that keeps\n# the current \"en-US\" run-mark language, followed by a new paragraph holding\n# the second line of text. Using InsertXML (instead of InsertParagraphAfter /\n# TypeText) lets us land this exact markup - no stray empty <w:r> and no\n# inherited formatting bleeding into the new text paragraph.\n$xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr></w:p><w:p><w:r><w:t>\u0414\u043e\u0431\u0430\u0432\u043b\u044f\u044e \u0432\u0442\u043e\u0440\u0443\u044e \u0441\u0442\u0440\u043e\u0447\u043a\u0443.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$r.InsertXML($xml)\n"}
